# Update countries & provincias Spain
# Refresh COVID-19 stats table and swap rank order of Irak / Camerun,
# plus bump the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Irak / Camerun (Irak now outranks Camerun) ---
# Row 68 used to be Camerun, row 69 used to be Irak; their names (and,
# correspondingly, the underlying shared-string order) flip along with
# the refreshed figures below.
$ws.Range("A68").Value = "Irak"
$ws.Range("A69").Value = "Camerun"

# --- Refreshed per-country figures (Casos totales, Nuevos casos, Casos
#     activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1676460
$ws.Range("C4").Value = 9632
$ws.Range("D4").Value = 451176
$ws.Range("E4").Value = 1126271
$ws.Range("G4").Value = 330
$ws.Range("H4").Value = 99013

# Row 10: Francia
$ws.Range("B10").Value = 182584
$ws.Range("C10").Value = 115
$ws.Range("D10").Value = 64617
$ws.Range("E10").Value = 89600
$ws.Range("G10").Value = 35
$ws.Range("H10").Value = 28367

# Row 11: Alemania
$ws.Range("B11").Value = 180167
$ws.Range("C11").Value = 181
$ws.Range("E11").Value = 11496

# Row 13: India
$ws.Range("B13").Value = 138536
$ws.Range("C13").Value = 7113
$ws.Range("D13").Value = 57692
$ws.Range("E13").Value = 76820

# Row 68: Irak (data refreshed in addition to the name swap above)
$ws.Range("B68").Value = 4469
$ws.Range("C68").Value = 197
$ws.Range("D68").Value = 2738
$ws.Range("E68").Value = 1571
$ws.Range("G68").Value = 8
$ws.Range("H68").Value = 160

# Row 69: Camerun
$ws.Range("B69").Value = 4400
$ws.Range("D69").Value = 1822
$ws.Range("E69").Value = 2419
$ws.Range("H69").Value = 159

# Row 76: Uzbekistan
$ws.Range("B76").Value = 3164
$ws.Range("C76").Value = 49
$ws.Range("E76").Value = 586

# Row 103: Sri Lanka
$ws.Range("B103").Value = 1140
$ws.Range("C103").Value = 51
$ws.Range("E103").Value = 457

# Row 151: Suazilandia
$ws.Range("B151").Value = 250
$ws.Range("C151").Value = 12
$ws.Range("D151").Value = 156
$ws.Range("E151").Value = 92

# --- Bump "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 20:35"
